$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them into numeric values
# (losing the literal "12.34"-style display and introducing FP rounding).
$numericLookingCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D17', 'D21', 'D22', 'D23', 'D24', 'D26', 'D28', 'D29', 'D31', 'D32', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.908.87'
$ws.Range('E2').Value = '  +4.07%  '
$ws.Range('D3').Value = '2.278.26'
$ws.Range('E3').Value = '  +4.47%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '251.43'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '0.638'
$ws.Range('E6').Value = '  +3.63%  '
$ws.Range('D7').Value = '71.82'
$ws.Range('E7').Value = '  +8.17%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.651'
$ws.Range('E9').Value = '  +13.54%  '
$ws.Range('D10').Value = '38.84'
$ws.Range('E10').Value = '  +6.92%  '
$ws.Range('D11').Value = '60.11'
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').Value = '0.0975'
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('D13').Value = '7.42'
$ws.Range('E13').Value = '  +7.58%  '
$ws.Range('D14').Value = '0.104'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').Value = '2.621.07'
$ws.Range('E15').Value = '  +4.63%  '
$ws.Range('D16').Value = '14.93'
$ws.Range('E16').Value = '  +3.87%  '
$ws.Range('D17').Value = '0.888'
$ws.Range('E17').Value = '  +4.12%  '
$ws.Range('D18').Value = '2.277.19'
$ws.Range('E18').Value = '  +5.48%  '
$ws.Range('D19').Value = '42.859.79'
$ws.Range('E19').Value = '  +4.09%  '
$ws.Range('E20').Value = '  +7.46%  '
$ws.Range('D21').Value = '6.31'
$ws.Range('E21').Value = '  +3.85%  '
$ws.Range('D22').Value = '73.48'
$ws.Range('E22').Value = '  +2.59%  '
$ws.Range('D23').Value = '236.42'
$ws.Range('E23').Value = '  +2.58%  '
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +5.53%  '
$ws.Range('E25').Value = '  +7.48%  '
$ws.Range('D26').Value = '11.52'
$ws.Range('E26').Value = '  +1.69%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '2.45'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').Value = '3.68'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('E30').Value = '  +5.54%  '
$ws.Range('D31').Value = '168.04'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '21.05'
$ws.Range('E32').Value = '  +3.98%  '
$ws.Range('D33').Value = '6.53'
$ws.Range('E33').Value = '  +13.28%  '
$ws.Range('E34').Value = '  +4.26%  '
$ws.Range('D35').Value = '31.53'
$ws.Range('E35').Value = '  +28.67%  '
$ws.Range('D36').Value = '0.0800'
$ws.Range('E36').Value = '  +8.63%  '
$ws.Range('D37').Value = '0.126'
$ws.Range('E37').Value = '  +3.86%  '
$ws.Range('D38').Value = '4.51'
$ws.Range('E38').Value = '  +14.05%  '
$ws.Range('D39').Value = '4.77'
$ws.Range('E39').Value = '  +5.41%  '
$ws.Range('E40').Value = '  +2.60%  '
$ws.Range('D41').Value = '13.39'
$ws.Range('E41').Value = '  +17.53%  '
$ws.Range('D42').Value = '2.32'
$ws.Range('E42').Value = '  +4.88%  '
$ws.Range('D43').Value = '5.84'
$ws.Range('E43').Value = '  +6.08%  '
$ws.Range('D44').Value = '0.210'
$ws.Range('E44').Value = '  +10.31%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '9.18'
$ws.Range('E45').Value = '  +7.62%  '
$ws.Range('B46').Value = 'MultiversX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D46').Value = '62.29'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').Value = '4.99'
$ws.Range('E47').Value = '  -8.03%  '
$ws.Range('E48').Value = '  +2.55%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('E50').Value = '  +3.93%  '
$ws.Range('E51').Value = '  +3.87%  '

# Restore default General number format / style now that the text is locked in,
# so the cells end up with no explicit style index (matching a plain text cell).
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
